$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15's date cell (B15) was the "last row" and used the date-only
# format; now that more rows follow, it should use the same
# date+time format as the rest of the column.
$ws.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New trade rows appended below the existing data (rows 16-25).
$newRows = @(
    @("FCL",       44972, 160,  256.9,    "Fineotex chemical limited"),
    @("LTTS",       44972, 11,   3622.04,  "L&t technology services limited"),
    @("ABBOTINDIA", 44970, -1,   20592.1,  "Abbott india limited"),
    @("ABBOTINDIA", 44970, -1,   20592.05, "Abbott india limited"),
    @("ABBOTINDIA", 44970, -2,   20591,    "Abbott india limited"),
    @("FCL",       44972, 5,    255.7,    "Fineotex chemical limited"),
    @("LTTS",       44972, -11,  3593.85,  "L&t technology services limited"),
    @("FCL",       44972, -165, 254,      "Fineotex chemical limited"),
    @("FCL",       44972, 160,  252.8,    "Fineotex chemical limited"),
    @("LTTS",       44972, 11,   3601.14,  "L&t technology services limited")
)

$startRow = 16
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# The new last row's date cell reverts to the date-only format,
# matching the original "last row" convention.
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("B" + $lastRow).NumberFormat = "YYYY-MM-DD"
